$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Coin / Link (text) columns: plain value assignment -------------------
$ws.Range("B6").Value = "KuCoinToken"
$ws.Range("C6").Value = "https://coinranking.com/coin/LOO6LmXd7G84Z+kucointoken-kcs"
$ws.Range("B7").Value = "GateToken"
$ws.Range("C7").Value = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
$ws.Range("B18").Value = "TigerCash"
$ws.Range("C18").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range("B19").Value = "LEO"
$ws.Range("C19").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("B20").Value = "BitpandaEcosystemToken"
$ws.Range("C20").Value = "https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best"
$ws.Range("B21").Value = "ProBitToken"
$ws.Range("C21").Value = "https://coinranking.com/coin/lQP4d6T2+probittoken-prob"
$ws.Range("B22").Value = "ZBToken"
$ws.Range("C22").Value = "https://coinranking.com/coin/CxmvOsCyENPso+zbtoken-zb"
$ws.Range("B23").Value = "CoinExToken"
$ws.Range("C23").Value = "https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet"

# --- Price / Volume(1h) columns: these are stored as plain text in the ----
# workbook (e.g. "39.20", "0.70%") even though they look numeric. Briefly
# force text format so Excel does not coerce them into numbers (which would
# drop significant trailing zeros / convert "%" strings), then restore the
# cell to its original General format so no stray formatting diff remains.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "331.03"
$ws.Range("D2").NumberFormat = "General"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "0.70%"
$ws.Range("E2").NumberFormat = "General"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "39.20"
$ws.Range("D3").NumberFormat = "General"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "-2.80%"
$ws.Range("E3").NumberFormat = "General"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "5.698"
$ws.Range("D4").NumberFormat = "General"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "1.84%"
$ws.Range("E4").NumberFormat = "General"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.08039"
$ws.Range("D5").NumberFormat = "General"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "-1.34%"
$ws.Range("E5").NumberFormat = "General"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "8.616"
$ws.Range("D6").NumberFormat = "General"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "-0.80%"
$ws.Range("E6").NumberFormat = "General"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "4.479"
$ws.Range("D7").NumberFormat = "General"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "-1.60%"
$ws.Range("E7").NumberFormat = "General"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "1.953"
$ws.Range("D8").NumberFormat = "General"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "-0.88%"
$ws.Range("E8").NumberFormat = "General"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2.971"
$ws.Range("D9").NumberFormat = "General"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "-0.23%"
$ws.Range("E9").NumberFormat = "General"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.9206"
$ws.Range("D10").NumberFormat = "General"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "-3.12%"
$ws.Range("E10").NumberFormat = "General"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.1237"
$ws.Range("D11").NumberFormat = "General"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "-2.64%"
$ws.Range("E11").NumberFormat = "General"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.1945"
$ws.Range("D12").NumberFormat = "General"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "-2.44%"
$ws.Range("E12").NumberFormat = "General"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "8.717"
$ws.Range("D13").NumberFormat = "General"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "16.18%"
$ws.Range("E13").NumberFormat = "General"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.09198"
$ws.Range("D14").NumberFormat = "General"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.03491"
$ws.Range("D15").NumberFormat = "General"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "-2.06%"
$ws.Range("E15").NumberFormat = "General"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.1050"
$ws.Range("D16").NumberFormat = "General"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "9.10%"
$ws.Range("E16").NumberFormat = "General"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.001295"
$ws.Range("D17").NumberFormat = "General"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "-1.32%"
$ws.Range("E17").NumberFormat = "General"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.006352"
$ws.Range("D18").NumberFormat = "General"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "4.42%"
$ws.Range("E18").NumberFormat = "General"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "3.363"
$ws.Range("D19").NumberFormat = "General"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "-0.38%"
$ws.Range("E19").NumberFormat = "General"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.3458"
$ws.Range("D20").NumberFormat = "General"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "-1.61%"
$ws.Range("E20").NumberFormat = "General"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.1371"
$ws.Range("D21").NumberFormat = "General"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "-1.93%"
$ws.Range("E21").NumberFormat = "General"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.2612"
$ws.Range("D22").NumberFormat = "General"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "5.20%"
$ws.Range("E22").NumberFormat = "General"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.04425"
$ws.Range("D23").NumberFormat = "General"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "-0.33%"
$ws.Range("E23").NumberFormat = "General"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.001252"
$ws.Range("D24").NumberFormat = "General"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "-0.06%"
$ws.Range("E24").NumberFormat = "General"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "4.51%"
$ws.Range("E25").NumberFormat = "General"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.0001200"
$ws.Range("D26").NumberFormat = "General"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "0.80%"
$ws.Range("E26").NumberFormat = "General"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.02573"
$ws.Range("D39").NumberFormat = "General"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "1.17%"
$ws.Range("E39").NumberFormat = "General"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.05441"
$ws.Range("D40").NumberFormat = "General"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "4.22%"
$ws.Range("E40").NumberFormat = "General"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.007528"
$ws.Range("D41").NumberFormat = "General"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "-3.89%"
$ws.Range("E41").NumberFormat = "General"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.009897"
$ws.Range("D42").NumberFormat = "General"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.1402"
$ws.Range("D43").NumberFormat = "General"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "-2.29%"
$ws.Range("E43").NumberFormat = "General"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.002110"
$ws.Range("D44").NumberFormat = "General"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "-3.69%"
$ws.Range("E44").NumberFormat = "General"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.01145"
$ws.Range("D45").NumberFormat = "General"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "15.44%"
$ws.Range("E45").NumberFormat = "General"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.00006810"
$ws.Range("D46").NumberFormat = "General"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "1.48%"
$ws.Range("E46").NumberFormat = "General"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "0.21%"
$ws.Range("E47").NumberFormat = "General"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.003055"
$ws.Range("D48").NumberFormat = "General"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "6.31%"
$ws.Range("E48").NumberFormat = "General"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.002281"
$ws.Range("D49").NumberFormat = "General"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "-0.82%"
$ws.Range("E49").NumberFormat = "General"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.00002102"
$ws.Range("D50").NumberFormat = "General"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "0.21%"
$ws.Range("E50").NumberFormat = "General"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0002002"
$ws.Range("D51").NumberFormat = "General"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "0.21%"
$ws.Range("E51").NumberFormat = "General"
